# Refresh the daily "cryptos" price/volume snapshot (GitHub Actions style
# update). Column D ("Price") cells are numeric-looking text (e.g. "0.637",
# "42.771.32" with thousands dots) that must stay TEXT, not auto-converted
# numbers, so we prefix with an apostrophe to force text entry and then
# reset the cell style back to "Normal" (undoing the quote-prefix
# formatting flag Excel would otherwise attach) so no stray style gets
# left on the cell. Column E ("Volume(1h)") values already contain
# padding spaces/percent signs so plain assignment keeps them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.771.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "'2.274.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.47%  "
$ws.Range("D5").Value = "'251.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Value = "'0.637"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("D7").Value = "'71.98"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.42%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.648"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +12.66%  "
$ws.Range("D10").Value = "'38.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.08%  "
$ws.Range("D11").Value = "'59.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  +4.20%  "
$ws.Range("D13").Value = "'7.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.06%  "
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "'2.614.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.32%  "
$ws.Range("D16").Value = "'14.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.43%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("D18").Value = "'2.272.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.11%  "
$ws.Range("D19").Value = "'42.703.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.90%  "
$ws.Range("D20").Value = "'0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.16%  "
$ws.Range("E21").Value = "  +3.57%  "
$ws.Range("D22").Value = "'73.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("D23").Value = "'234.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.27%  "
$ws.Range("D25").Value = "'3.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.15%  "
$ws.Range("D26").Value = "'11.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "'2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.73%  "
$ws.Range("D31").Value = "'168.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'21.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").Value = "'6.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.18%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.93%  "
$ws.Range("D35").Value = "'0.0801"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.13%  "
$ws.Range("D36").Value = "'30.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +26.17%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").Value = "'4.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.28%  "
$ws.Range("E39").Value = "  +5.05%  "
$ws.Range("E40").Value = "  +1.90%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'13.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.02%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "'2.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.00%  "
$ws.Range("E43").Value = "  +6.08%  "
$ws.Range("D44").Value = "'0.212"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.55%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.69%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "'4.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.47%  "
$ws.Range("D47").Value = "'61.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  +2.86%  "
$ws.Range("D49").Value = "'1.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +3.90%  "
